# "Changes of Config file"
# RTE_STG.xlsx: the RTE tracking-number config rows on "RTECreation" and
# "SearchRTE" are refreshed with a newly generated set of IDs. We write the
# new values as text-formula results and paste them back as values so the
# cells land as plain text (shared strings) just like the originals, without
# disturbing their (default) cell formatting/style.

$wb  = $excel.ActiveWorkbook
$xlPasteValues = -4163

$wsCreation = $wb.Worksheets.Item("RTECreation")
$wsSearch   = $wb.Worksheets.Item("SearchRTE")

# RTECreation!C2:C3 -> new RWTrackingNo values
$wsCreation.Range("C2").Formula = "=""126040034"""
$wsCreation.Range("C3").Formula = "=""126040056"""
$wsCreation.Range("C2:C3").Copy()
$wsCreation.Range("C2:C3").PasteSpecial($xlPasteValues)

# SearchRTE!A2:D3 -> matching new RWTrackingNo/JobID/PickUPID/BOLNo row set
$wsSearch.Range("A2").Formula = "=""126040034"""
$wsSearch.Range("B2").Formula = "=""32399896"""
$wsSearch.Range("C2").Formula = "=""3405134"""
$wsSearch.Range("D2").Formula = "=""126040045"""
$wsSearch.Range("A3").Formula = "=""126040056"""
$wsSearch.Range("B3").Formula = "=""32399897"""
$wsSearch.Range("C3").Formula = "=""3405135"""
$wsSearch.Range("D3").Formula = "=""126040067"""
$wsSearch.Range("A2:D3").Copy()
$wsSearch.Range("A2:D3").PasteSpecial($xlPasteValues)

$excel.CutCopyMode = $false
